# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The detail table (rows 16-42, cols B:J) is re-sorted: instead of being
# grouped by worker (period descending within each worker), it is now
# grouped by period (ascending), and within each period by worker
# (Eduardo Jose Torres Porto, William Ricardo Morales Lopez, Caleb
# Miranda Avila) -- except period 1907/1908/1909/1910 which only have
# the rows that originally existed for that period. Column H:J stay
# blank; F/G values travel together with the (worker, period) pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$data = @(
    ,@("CC","1050961099","EDUARDO JOSE TORRES PORTO","1907",17667,828116)
    ,@("CC","1001833583","WILLIAM RICARDO MORALES LOPEZ","1907",17667,828116)
    ,@("CC","1002180038","CALEB MIRANDA AVILA","1907",17667,828116)
    ,@("CC","1002180038","CALEB MIRANDA AVILA","1908",33125,828116)
    ,@("CC","1002180038","CALEB MIRANDA AVILA","1909",33125,828116)
    ,@("CC","1002180038","CALEB MIRANDA AVILA","1910",33125,828116)
    ,@("CC","1050961099","EDUARDO JOSE TORRES PORTO","1911",33125,828116)
    ,@("CC","1001833583","WILLIAM RICARDO MORALES LOPEZ","1911",33125,828116)
    ,@("CC","1002180038","CALEB MIRANDA AVILA","1911",33125,828116)
    ,@("CC","1050961099","EDUARDO JOSE TORRES PORTO","1912",33125,828116)
    ,@("CC","1001833583","WILLIAM RICARDO MORALES LOPEZ","1912",33125,828116)
    ,@("CC","1002180038","CALEB MIRANDA AVILA","1912",33125,828116)
    ,@("CC","1050961099","EDUARDO JOSE TORRES PORTO","2001",33125,828116)
    ,@("CC","1001833583","WILLIAM RICARDO MORALES LOPEZ","2001",33125,828116)
    ,@("CC","1002180038","CALEB MIRANDA AVILA","2001",33125,828116)
    ,@("CC","1050961099","EDUARDO JOSE TORRES PORTO","2002",33125,828116)
    ,@("CC","1001833583","WILLIAM RICARDO MORALES LOPEZ","2002",33125,828116)
    ,@("CC","1002180038","CALEB MIRANDA AVILA","2002",33125,828116)
    ,@("CC","1050961099","EDUARDO JOSE TORRES PORTO","2003",33125,828116)
    ,@("CC","1001833583","WILLIAM RICARDO MORALES LOPEZ","2003",33125,828116)
    ,@("CC","1002180038","CALEB MIRANDA AVILA","2003",33125,828116)
    ,@("CC","1050961099","EDUARDO JOSE TORRES PORTO","2004",33125,828116)
    ,@("CC","1001833583","WILLIAM RICARDO MORALES LOPEZ","2004",33125,828116)
    ,@("CC","1002180038","CALEB MIRANDA AVILA","2004",33125,828116)
    ,@("CC","1050961099","EDUARDO JOSE TORRES PORTO","2005",20979,828116)
    ,@("CC","1001833583","WILLIAM RICARDO MORALES LOPEZ","2005",20979,828116)
    ,@("CC","1002180038","CALEB MIRANDA AVILA","2005",20979,828116)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 2).Value = $rec[0]   # B: Tipo Doc Trabajador
    $ws.Cells.Item($row, 3).Value = $rec[1]   # C: N Doc Trabajador
    $ws.Cells.Item($row, 4).Value = $rec[2]   # D: Nombre Trabajador
    $ws.Cells.Item($row, 5).Value = $rec[3]   # E: Periodo Mora
    $ws.Cells.Item($row, 6).Value = $rec[4]   # F: Valor Mora
    $ws.Cells.Item($row, 7).Value = $rec[5]   # G: Salario Basico
}
